# Weekly update: insert a new price record (Ají - Cristal) as a new row 103,
# pushing the existing rows 103..210 down to 104..211.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above the current row 103; this shifts rows
# 103..210 down to 104..211 (matching the widened A1:R211 dimension).
$ws.Range("A103").EntireRow.Insert()

# Populate the newly inserted row 103 with the new record's data.
$ws.Range("A103").Value = 5
$ws.Range("B103").Value = "Macroferia Regional de Talca"
$ws.Range("C103").Value = "Maule"
$ws.Range("D103").Value = 44629
$ws.Range("E103").Value = 7
$ws.Range("F103").Value = 100112021
$ws.Range("G103").Value = "Ají"
$ws.Range("H103").Value = "Cristal"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 150
$ws.Range("K103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 15000
$ws.Range("N103").Value = "`$/saco 25 kilos"
$ws.Range("O103").Value = "Región del Maule"
$ws.Range("P103").Value = 600
$ws.Range("Q103").Value = 25
$ws.Range("R103").Value = "Hortaliza"
